$d = $word.ActiveDocument

# 1. ACTA Nº: "1" -> "12"  (paragraph 1)
$p = $d.Paragraphs(1)
$p.Range.Find.Execute("1", $true, $true, $false, $false, $false, $true, 1, $false, "12", 2)

# 2. FECHA: "12/07/2023" -> "28/08/2023"  (paragraph 4)
$p = $d.Paragraphs(4)
$p.Range.Find.Execute("12/07/2023", $true, $false, $false, $false, $false, $true, 1, $false, "28/08/2023", 2)

# 3. CONVOCA: "admin" -> "Mario"  (paragraph 5)
$p = $d.Paragraphs(5)
$p.Range.Find.Execute("admin", $true, $false, $false, $false, $false, $true, 1, $false, "Mario", 2)

# 4. LUGAR: "Envigado" -> "envigado"  (paragraph 7)
$p = $d.Paragraphs(7)
$p.Range.Find.Execute("Envigado", $true, $false, $false, $false, $false, $true, 1, $false, "envigado", 2)

# 5. HORA INICIO: "07:17:17" -> "03:13:48"  (paragraph 8)
$p = $d.Paragraphs(8)
$p.Range.Find.Execute("07:17:17", $true, $false, $false, $false, $false, $true, 1, $false, "03:13:48", 2)

# 6. "sa" -> "prueba"  (paragraph 11, ASISTENTES entry)
$p = $d.Paragraphs(11)
$p.Range.Find.Execute("sa", $true, $false, $false, $false, $false, $true, 1, $false, "prueba", 2)

# 7. "asa" -> "dsdsdsds"  (paragraph 15, ORDEN DEL DIA entry)
$p = $d.Paragraphs(15)
$p.Range.Find.Execute("asa", $true, $false, $false, $false, $false, $true, 1, $false, "dsdsdsds", 2)

# 8. "asa" -> "sdsdsdsds"  (paragraph 22, DESARROLLO ORDEN DEL DIA entry)
$p = $d.Paragraphs(22)
$p.Range.Find.Execute("asa", $true, $false, $false, $false, $false, $true, 1, $false, "sdsdsdsds", 2)
